$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 478   # 合肥·国乙only宇宙心动（含夜场）: 473 -> 478
$wsExpo.Range("F4").Value = 14    # 合肥·IE动漫嘉年华: 13 -> 14
$wsExpo.Range("F9").Value = 444   # 合肥·首届AT次元时代动漫游戏嘉年华: 322 -> 444

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 85    # 合肥·CrossingX意次元｜乐队番ONLY同人: 84 -> 85

# Sheet "全部类型" (All types) - aggregated view of the above sheets
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 85     # 合肥·CrossingX意次元｜乐队番ONLY同人: 84 -> 85
$wsAll.Range("F4").Value = 478    # 合肥·国乙only宇宙心动（含夜场）: 473 -> 478
$wsAll.Range("F5").Value = 14     # 合肥·IE动漫嘉年华: 13 -> 14
$wsAll.Range("F10").Value = 444   # 合肥·首届AT次元时代动漫游戏嘉年华: 322 -> 444
